# "Generate Report for handoff"
# Updates the localization-status report so that row 3 (b.md.md) reflects
# that it is now ready for handoff, instead of the stale
# "Handed back: in sync with en-US" status that was pointing at the wrong
# (a.md) handoff artifact.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: update the Status columns (zh-cn / de-de) for b.md.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet: update Status / Latest Handoff File / Latest Handoff
# Datetime for the b.md.md row (row 3), and repoint the hyperlink
# display text on the handoff-file cell to the new b.md artifact.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-18 12:21:00"

foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: same set of updates, localized for the de-de artifact.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-18 12:21:10"

foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
